$d = $word.ActiveDocument

$p3 = $d.Paragraphs(3)

# -----------------------------------------------------------------
# First, insert the two new list items right after paragraph 3 while
# its text is still simple (un-split). This keeps their runs free of
# any leftover direct-formatting markers picked up from edits made
# further up in the paragraph.
# -----------------------------------------------------------------
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "once package.json has been updated with all the dependencies, when running the repo on some other machine, just do npm install, which will search package.json for dependencies and install them all"

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "Note that a single repo can have different package.json files and node_modules in different folders, as required. "

# -----------------------------------------------------------------
# Now split "add the node modules to gitignore" into three runs and
# extend the sentence, matching:
#      "add the node_" | "modules to gitignore" |
#      " so that node_modules is not updated on git when committing"
# -----------------------------------------------------------------
$p3start = $p3.Range.Start

# Turn the space between "node" and "modules" into an underscore so the
# paragraph text becomes "add the node_modules to gitignore".
$space = $d.Range($p3start + 12, $p3start + 13)
$space.Text = "_"

# Append the extra clause at the end of the paragraph (before the pilcrow).
$p3TextEnd = $p3.Range.End - 1
$tail = $d.Range($p3TextEnd, $p3TextEnd)
$tail.InsertAfter(" so that node_modules is not updated on git when committing")

# Force run boundaries at the two seams. Toggling a character attribute
# on a sub-range and then clearing it again keeps the text/formatting
# identical to the surrounding runs while still splitting them apart,
# exactly like Word leaves behind separate <w:r> elements after an
# editing session that touched those spans.
$split1 = $d.Range($p3start + 13, $p3start + 33)
$split1.Bold = 1
$split1.Bold = 0

$split2 = $d.Range($p3start + 33, $p3.Range.End - 1)
$split2.Bold = 1
$split2.Bold = 0

# -----------------------------------------------------------------
# Split "Note that a single repo can have ..." into four runs:
# "Note that a single " | "r" | "epo" |
# " can have different package.json files and node_modules in
# different folders, as required. "
# -----------------------------------------------------------------
$p5start = $p5.Range.Start

$r1 = $d.Range($p5start + 19, $p5start + 20)
$r1.Bold = 1
$r1.Bold = 0

$r2 = $d.Range($p5start + 20, $p5start + 23)
$r2.Bold = 1
$r2.Bold = 0

$r3 = $d.Range($p5start + 23, $p5.Range.End - 1)
$r3.Bold = 1
$r3.Bold = 0

Write-Output "done"
